# Apply updated cryptocurrency price/volume data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "6.40"); force
# text format first so Excel does not silently coerce it to a Number
# and drop significant trailing zeros.
$dRefs = @("D2","D3","D5","D6","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D25","D27","D30","D31","D32","D35","D38","D40","D41","D42","D43","D44","D45","D47","D48","D50","D51")
foreach ($ref in $dRefs) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = '44.552.92'
$ws.Range("E2").Value = '  +3.76%  '
$ws.Range("D3").Value = '2.422.35'
$ws.Range("E3").Value = '  +2.55%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '314.96'
$ws.Range("E5").Value = '  +4.14%  '
$ws.Range("D6").Value = '101.01'
$ws.Range("E6").Value = '  +5.67%  '
$ws.Range("E7").Value = '  +2.26%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '0.532'
$ws.Range("E9").Value = '  +10.18%  '
$ws.Range("D10").Value = '35.35'
$ws.Range("E10").Value = '  +3.87%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.0799'
$ws.Range("E11").Value = '  +1.89%  '
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").Value = '19.05'
$ws.Range("E12").Value = '  +3.33%  '
$ws.Range("D13").Value = '0.121'
$ws.Range("E13").Value = '  -2.35%  '
$ws.Range("D14").Value = '6.97'
$ws.Range("E14").Value = '  +4.23%  '
$ws.Range("D15").Value = '2.799.83'
$ws.Range("E15").Value = '  +2.62%  '
$ws.Range("D16").Value = '2.421.99'
$ws.Range("E16").Value = '  +2.66%  '
$ws.Range("D17").Value = '0.831'
$ws.Range("E17").Value = '  +5.33%  '
$ws.Range("D18").Value = '44.440.41'
$ws.Range("E18").Value = '  +3.62%  '
$ws.Range("D19").Value = '12.45'
$ws.Range("E19").Value = '  +4.44%  '
$ws.Range("D20").Value = '6.40'
$ws.Range("E20").Value = '  +2.27%  '
$ws.Range("D21").Value = '0.0₃0917'
$ws.Range("E21").Value = '  +3.72%  '
$ws.Range("D22").Value = '68.79'
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").Value = '242.30'
$ws.Range("E23").Value = '  +3.12%  '
$ws.Range("E24").Value = '  +5.30%  '
$ws.Range("D25").Value = '2.48'
$ws.Range("E25").Value = '  +2.09%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '25.20'
$ws.Range("E27").Value = '  +3.21%  '
$ws.Range("E28").Value = '  -3.87%  '
$ws.Range("E29").Value = '  +2.77%  '
$ws.Range("D30").Value = '33.17'
$ws.Range("E30").Value = '  +4.03%  '
$ws.Range("D31").Value = '48.40'
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("D32").Value = '0.124'
$ws.Range("E32").Value = '  +19.95%  '
$ws.Range("E33").Value = '  +9.84%  '
$ws.Range("E34").Value = '  +3.42%  '
$ws.Range("D35").Value = '0.0774'
$ws.Range("E35").Value = '  +7.80%  '
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("E37").Value = '  +2.66%  '
$ws.Range("D38").Value = '4.48'
$ws.Range("E38").Value = '  +4.12%  '
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("D40").Value = '2.22'
$ws.Range("E40").Value = '  -2.61%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.109'
$ws.Range("E41").Value = '  +1.72%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '118.77'
$ws.Range("E42").Value = '  -7.01%  '
$ws.Range("D43").Value = '20.96'
$ws.Range("E43").Value = '  -0.80%  '
$ws.Range("D44").Value = '0.0290'
$ws.Range("E44").Value = '  +4.37%  '
$ws.Range("D45").Value = '1.943.97'
$ws.Range("E45").Value = '  +0.92%  '
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("D47").Value = '2.95'
$ws.Range("E47").Value = '  +9.28%  '
$ws.Range("D48").Value = '9.43'
$ws.Range("E48").Value = '  +2.67%  '
$ws.Range("E49").Value = '  +10.04%  '
$ws.Range("D50").Value = '54.76'
$ws.Range("E50").Value = '  +6.64%  '
$ws.Range("D51").Value = '74.68'
$ws.Range("E51").Value = '  +4.54%  '
